# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "59.099.37"
Set-TextValue "E2" "  -0.33%  "
Set-TextValue "D3" "2.522.09"
Set-TextValue "E3" "  +0.28%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "536.72"
Set-TextValue "E5" "  +0.30%  "
Set-TextValue "D6" "137.61"
Set-TextValue "E6" "  -1.30%  "
Set-TextValue "E7" "  +0.15%  "
Set-TextValue "D8" "0.568"
Set-TextValue "E8" "  +0.49%  "
Set-TextValue "D9" "2.520.59"
Set-TextValue "E9" "  +0.09%  "
Set-TextValue "E10" "  +0.22%  "
Set-TextValue "D11" "0.157"
Set-TextValue "E11" "  -2.50%  "
Set-TextValue "E12" "  -1.45%  "
Set-TextValue "E13" "  -2.35%  "
Set-TextValue "D14" "2.970.41"
Set-TextValue "E14" "  +0.28%  "
Set-TextValue "D15" "23.01"
Set-TextValue "E15" "  -1.99%  "
Set-TextValue "D16" "59.052.43"
Set-TextValue "E16" "  -0.25%  "
Set-TextValue "E17" "  -1.41%  "
Set-TextValue "D18" "2.534.95"
Set-TextValue "E18" "  +0.72%  "
Set-TextValue "D19" "11.14"
Set-TextValue "E19" "  +0.00%  "
Set-TextValue "E20" "  -0.45%  "
Set-TextValue "D21" "324.12"
Set-TextValue "E21" "  -0.33%  "
Set-TextValue "E22" "  +0.02%  "
Set-TextValue "D23" "5.94"
Set-TextValue "E23" "  +2.49%  "
Set-TextValue "D24" "65.66"
Set-TextValue "E24" "  +2.96%  "
Set-TextValue "E25" "  -1.49%  "
Set-TextValue "E26" "  -0.09%  "
Set-TextValue "E27" "  +0.15%  "
Set-TextValue "D28" "7.61"
Set-TextValue "E28" "  -3.00%  "
Set-TextValue "E29" "  -2.65%  "
Set-TextValue "E30" "  -1.13%  "
Set-TextValue "E31" "  -0.40%  "
Set-TextValue "B32" "Monero"
Set-TextValue "C32" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "164.98"
Set-TextValue "E32" "  -0.08%  "
Set-TextValue "B33" "Fetch.AI"
Set-TextValue "C33" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D33" "1.17"
Set-TextValue "E33" "  +5.21%  "
Set-TextValue "D34" "1.48"
Set-TextValue "E34" "  +1.30%  "
Set-TextValue "D35" "0.999"
Set-TextValue "E35" "  +0.04%  "
Set-TextValue "D36" "18.45"
Set-TextValue "E36" "  -0.38%  "
Set-TextValue "D37" "4.10"
Set-TextValue "E37" "  -3.61%  "
Set-TextValue "E38" "  -2.49%  "
Set-TextValue "D39" "36.73"
Set-TextValue "E39" "  -0.44%  "
Set-TextValue "D40" "0.818"
Set-TextValue "E40" "  +0.11%  "
Set-TextValue "D41" "3.63"
Set-TextValue "E41" "  -1.81%  "
Set-TextValue "D42" "286.34"
Set-TextValue "E42" "  +2.51%  "
Set-TextValue "D43" "5.20"
Set-TextValue "E43" "  -1.11%  "
Set-TextValue "D44" "132.32"
Set-TextValue "E44" "  +7.59%  "
Set-TextValue "D45" "0.999"
Set-TextValue "E45" "  +0.15%  "
Set-TextValue "D46" "0.605"
Set-TextValue "E46" "  +1.23%  "
Set-TextValue "D47" "10.87"
Set-TextValue "E47" "  +0.00%  "
Set-TextValue "E48" "  -0.32%  "
Set-TextValue "E49" "  -0.88%  "
Set-TextValue "D50" "0.0221"
Set-TextValue "E50" "  -1.52%  "
Set-TextValue "D51" "17.27"
Set-TextValue "E51" "  -3.04%  "
